$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.73
$ws.Range("I2").Value = 4.33
$ws.Range("J2").Value = 2.3
$ws.Range("L2").Value = 4.5
$ws.Range("U2").Value = 1.67
$ws.Range("V2").Value = 2.1
$ws.Range("X2").Value = 9
$ws.Range("AK2").Value = 51
$ws.Range("AL2").Value = 34
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 23
